# "Add files via upload" -- reproduce the edits made to the two body
# paragraphs of "KK Ethical reflection.docx":
#
#   1. Strip the direct/legacy formatting (NormalWeb paragraph style and an
#      explicit color/sz/szCs run-properties block repeated on every run) so
#      the paragraphs and runs fall back to the document defaults.
#   2. Extend the final sentence of paragraph 2 with two new sentences about
#      nanoparticles crossing the blood-brain barrier.
#   3. Append a further sentence about sunscreen as new runs (with the usual
#      gramStart/gramEnd proofing marks Word puts around "The").
#   4. Insert one extra empty paragraph before the empty paragraph that was
#      already at the end of the document.
#
# Find/Replace can append text onto an existing run (keeping its rPr), but it
# cannot *remove* an existing <w:rPr>/<w:pPr>, and it cannot create new runs
# that carry no run properties at all. The faithful way to reproduce that
# exact run layout is to hand Word the replacement WordprocessingML for the
# whole body via Range.InsertXML (a regular member of the Word object model).

$d = $word.ActiveDocument

$para1 = '<w:p w14:paraId="684813E0" w14:textId="77777777" w:rsidR="000D23CD" w:rsidRDefault="000D23CD" w:rsidP="000D23CD"><w:r><w:t xml:space="preserve">During this project it has given me the opportunity to do some research on what some of the risks could be in nanotechnology in medicine. </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>I''ve</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> never heard of Nanotechnology and was a bit intimidated by the name once we chose this as our group topic, however, I was determined to explore what the world of nanotechnology had to teach me.</w:t></w:r></w:p>'
$para2 = '<w:p w14:paraId="302F6001" w14:textId="77777777" w:rsidR="000D23CD" w:rsidRDefault="000D23CD" w:rsidP="000D23CD"><w:r><w:t xml:space="preserve">I found that an ethical issue in nanomedicine is protecting </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>patients</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> private information where medical records of body organs are stored electronically of the patients results. In addition, it has been revealed that when clinicians have a consultation remotely with their patients to discuss test results or diagnosis that has recently been discovered, this data is already stored in the system which is allocated against the </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>patients</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> file. The quantity of the storage is quite </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>large,</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> therefore it is recommended to have a highly sourced secure system to protect patient information from a cyber-attack. Nanoparticles can be a risk when it enters the body as this could attack an important body organ which is the brain through a new drug. If the new drug passed through the blood brain barrier which protects toxins trying to enter the brain, this will affect the spinal cord or lead to another condition.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>The</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> thought that nanoparticles can enter the body through using sunscreen is something to be wary of.</w:t></w:r></w:p>'
$trailingParas = '<w:p/><w:p w14:paraId="3B399186" w14:textId="77777777" w:rsidR="00050467" w:rsidRDefault="004D14B0"/><w:p/>'

$newBody = $para1 + $para2 + $trailingParas

$packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>{0}</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$packageXml = $packageXml -replace [regex]::Escape("{0}"), $newBody

# Replacing the whole body keeps paragraph/run formatting exactly as built
# above (no stray empty <w:rPr/> placeholders), and Word keeps the existing
# section properties (headers, page size, ...) that already followed the
# last paragraph of the body.
$d.Content.InsertXML($packageXml)
